# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2: Home) with new simulation totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 232
$wsOff.Range("C2").Value = 145
$wsOff.Range("D2").Value = 47
$wsOff.Range("E2").Value = 24

# Update DEF sheet (row 2: Home) with new simulation totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 200
$wsDef.Range("C2").Value = 150
$wsDef.Range("D2").Value = 51
$wsDef.Range("E2").Value = 33
